# Updates currentAveragePrice / Leve price & profit columns (H..N) across the
# ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets to reflect refreshed market-board data
# pulled by the scheduled runner. Each row is located by its unique "Leve
# Item ID" value in column G so edits land on the correct leve even though
# row numbers repeat across sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4: Leve Item ID (G4) = 5470
if ($ws.Cells.Item(4, 7).Value2 -eq 5470) {
    $ws.Range("H4").Value = 790.5
    $ws.Range("I4").Value = 189.14285
    $ws.Range("K4").Value = 189.14285
    $ws.Range("M4").Value = -75.14285000000001
}
# Row 9: Leve Item ID (G9) = 5487
if ($ws.Cells.Item(9, 7).Value2 -eq 5487) {
    $ws.Range("H9").Value = 935.63635
    $ws.Range("I9").Value = 161.25
    $ws.Range("K9").Value = 161.25
    $ws.Range("M9").Value = 7.75
}
# Row 17: Leve Item ID (G17) = 38956
if ($ws.Cells.Item(17, 7).Value2 -eq 38956) {
    $ws.Range("H17").Value = 920.48
    $ws.Range("J17").Value = 920.48
    $ws.Range("L17").Value = 2761.44
    $ws.Range("N17").Value = -3097.44
}
# Row 40: Leve Item ID (G40) = 5505
if ($ws.Cells.Item(40, 7).Value2 -eq 5505) {
    $ws.Range("H40").Value = 2637.8696
    $ws.Range("I40").Value = 1745.5
    $ws.Range("J40").Value = 2825.7368
    $ws.Range("K40").Value = 1745.5
    $ws.Range("L40").Value = 2825.7368
    $ws.Range("M40").Value = -1570.5
    $ws.Range("N40").Value = -3175.7368
}
# Row 41: Leve Item ID (G41) = 5478
if ($ws.Cells.Item(41, 7).Value2 -eq 5478) {
    $ws.Range("H41").Value = 482.5
    $ws.Range("I41").Value = 248.1
    $ws.Range("J41").Value = 1068.5
    $ws.Range("K41").Value = 248.1
    $ws.Range("L41").Value = 1068.5
    $ws.Range("M41").Value = 191.9
    $ws.Range("N41").Value = -1948.5
}
# Row 64: Leve Item ID (G64) = 5506
if ($ws.Cells.Item(64, 7).Value2 -eq 5506) {
    $ws.Range("H64").Value = 4875
    $ws.Range("I64").Value = 4525
    $ws.Range("J64").Value = 5750
    $ws.Range("K64").Value = 4525
    $ws.Range("L64").Value = 5750
    $ws.Range("M64").Value = -4277
    $ws.Range("N64").Value = -6246
}
# Row 67: Leve Item ID (G67) = 5506
if ($ws.Cells.Item(67, 7).Value2 -eq 5506) {
    $ws.Range("H67").Value = 4875
    $ws.Range("I67").Value = 4525
    $ws.Range("J67").Value = 5750
    $ws.Range("K67").Value = 4525
    $ws.Range("L67").Value = 5750
    $ws.Range("M67").Value = -3667
    $ws.Range("N67").Value = -7466
}
# Row 100: Leve Item ID (G100) = 19906
if ($ws.Cells.Item(100, 7).Value2 -eq 19906) {
    $ws.Range("H100").Value = 5888.6665
    $ws.Range("J100").Value = 5888.6665
    $ws.Range("L100").Value = 5888.6665
    $ws.Range("N100").Value = -6970.6665
}
# Row 111: Leve Item ID (G111) = 27768
if ($ws.Cells.Item(111, 7).Value2 -eq 27768) {
    $ws.Range("H111").Value = 4268.524
    $ws.Range("I111").Value = 4281.5
    $ws.Range("J111").Value = 4260.5386
    $ws.Range("K111").Value = 12844.5
    $ws.Range("L111").Value = 12781.6158
    $ws.Range("M111").Value = -9777.5
    $ws.Range("N111").Value = -18915.6158
}
# Row 113: Leve Item ID (G113) = 27775
if ($ws.Cells.Item(113, 7).Value2 -eq 27775) {
    $ws.Range("H113").Value = 6973.5186
    $ws.Range("I113").Value = 6064.3335
    $ws.Range("K113").Value = 6064.3335
    $ws.Range("M113").Value = -2810.3335
}
# Row 124: Leve Item ID (G124) = 34241
if ($ws.Cells.Item(124, 7).Value2 -eq 34241) {
    $ws.Range("H124").Value = 83920
    $ws.Range("J124").Value = 83920
    $ws.Range("L124").Value = 83920
    $ws.Range("N124").Value = -93740
}
# Row 138: Leve Item ID (G138) = 44169
if ($ws.Cells.Item(138, 7).Value2 -eq 44169) {
    $ws.Range("H138").Value = 2963.6304
    $ws.Range("J138").Value = 4505.7334
    $ws.Range("L138").Value = 13517.2002
    $ws.Range("N138").Value = -23797.2002
}

$ws = $wb.Worksheets.Item("ARM")
# Row 17: Leve Item ID (G17) = 2495
if ($ws.Cells.Item(17, 7).Value2 -eq 2495) {
    $ws.Range("H17").Value = 300
    $ws.Range("I17").Value = 300
    $ws.Range("J17").Value = 0
    $ws.Range("K17").Value = 300
    $ws.Range("L17").Value = 0
    $ws.Range("M17").Value = -127
    $ws.Range("N17").ClearContents()
}
# Row 18: Leve Item ID (G18) = 2470
if ($ws.Cells.Item(18, 7).Value2 -eq 2470) {
    $ws.Range("H18").Value = 749.5
    $ws.Range("J18").Value = 749.5
    $ws.Range("L18").Value = 749.5
    $ws.Range("N18").Value = -1393.5
}
# Row 45: Leve Item ID (G45) = 27714
if ($ws.Cells.Item(45, 7).Value2 -eq 27714) {
    $ws.Range("H45").Value = 2130.25
    $ws.Range("J45").Value = 2719.5
    $ws.Range("L45").Value = 2719.5
    $ws.Range("N45").Value = -3473.5
}

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Leve Item ID (G22) = 5092
if ($ws.Cells.Item(22, 7).Value2 -eq 5092) {
    $ws.Range("H22").Value = 712.5
    $ws.Range("I22").Value = 533.3333
    $ws.Range("J22").Value = 1250
    $ws.Range("K22").Value = 533.3333
    $ws.Range("L22").Value = 1250
    $ws.Range("M22").Value = -360.3333
    $ws.Range("N22").Value = -1596
}
# Row 98: Leve Item ID (G98) = 19545
if ($ws.Cells.Item(98, 7).Value2 -eq 19545) {
    $ws.Range("H98").Value = 66540
    $ws.Range("J98").Value = 66540
    $ws.Range("L98").Value = 66540
    $ws.Range("N98").Value = -72530
}

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Leve Item ID (G68) = 12895
if ($ws.Cells.Item(68, 7).Value2 -eq 12895) {
    $ws.Range("H68").Value = 4614.4287
    $ws.Range("I68").Value = 876
    $ws.Range("J68").Value = 6109.8
    $ws.Range("K68").Value = 2628
    $ws.Range("L68").Value = 18329.4
    $ws.Range("M68").Value = -1817
    $ws.Range("N68").Value = -19951.4
}
# Row 71: Leve Item ID (G71) = 12895
if ($ws.Cells.Item(71, 7).Value2 -eq 12895) {
    $ws.Range("H71").Value = 4614.4287
    $ws.Range("I71").Value = 876
    $ws.Range("J71").Value = 6109.8
    $ws.Range("K71").Value = 7884
    $ws.Range("L71").Value = 54988.2
    $ws.Range("M71").Value = -3828
    $ws.Range("N71").Value = -63100.2
}
# Row 118: Leve Item ID (G118) = 27872
if ($ws.Cells.Item(118, 7).Value2 -eq 27872) {
    $ws.Range("H118").Value = 5148.25
    $ws.Range("I118").Value = 3975
    $ws.Range("K118").Value = 11925
    $ws.Range("M118").Value = -10682
}

$ws = $wb.Worksheets.Item("GSM")
# Row 97: Leve Item ID (G97) = 19940
if ($ws.Cells.Item(97, 7).Value2 -eq 19940) {
    $ws.Range("H97").Value = 1115.2
    $ws.Range("I97").Value = 464.57144
    $ws.Range("J97").Value = 2633.3333
    $ws.Range("K97").Value = 464.57144
    $ws.Range("L97").Value = 2633.3333
    $ws.Range("M97").Value = 31.42856
    $ws.Range("N97").Value = -3625.3333
}

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Leve Item ID (G22) = 5277
if ($ws.Cells.Item(22, 7).Value2 -eq 5277) {
    $ws.Range("H22").Value = 1340.4348
    $ws.Range("I22").Value = 413.5
    $ws.Range("J22").Value = 1667.5883
    $ws.Range("K22").Value = 413.5
    $ws.Range("L22").Value = 1667.5883
    $ws.Range("M22").Value = -118.5
    $ws.Range("N22").Value = -2257.5883
}
# Row 27: Leve Item ID (G27) = 5277
if ($ws.Cells.Item(27, 7).Value2 -eq 5277) {
    $ws.Range("H27").Value = 1340.4348
    $ws.Range("I27").Value = 413.5
    $ws.Range("J27").Value = 1667.5883
    $ws.Range("K27").Value = 413.5
    $ws.Range("L27").Value = 1667.5883
    $ws.Range("M27").Value = -306.5
    $ws.Range("N27").Value = -1881.5883
}
# Row 36: Leve Item ID (G36) = 34261
if ($ws.Cells.Item(36, 7).Value2 -eq 34261) {
    $ws.Range("H36").Value = 69582.664
    $ws.Range("J36").Value = 69582.664
    $ws.Range("L36").Value = 69582.664
    $ws.Range("N36").Value = -70706.664
}
# Row 46: Leve Item ID (G46) = 5282
if ($ws.Cells.Item(46, 7).Value2 -eq 5282) {
    $ws.Range("H46").Value = 2371.4211
    $ws.Range("I46").Value = 1133.3334
    $ws.Range("J46").Value = 2603.5625
    $ws.Range("K46").Value = 1133.3334
    $ws.Range("L46").Value = 2603.5625
    $ws.Range("M46").Value = -945.3334
    $ws.Range("N46").Value = -2979.5625
}
# Row 122: Leve Item ID (G122) = 36247
if ($ws.Cells.Item(122, 7).Value2 -eq 36247) {
    $ws.Range("H122").Value = 6444
    $ws.Range("I122").Value = 6283
    $ws.Range("J122").Value = 6478.5
    $ws.Range("K122").Value = 18849
    $ws.Range("L122").Value = 19435.5
    $ws.Range("M122").Value = -16399
    $ws.Range("N122").Value = -24335.5
}
# Row 132: Leve Item ID (G132) = 44058
if ($ws.Cells.Item(132, 7).Value2 -eq 44058) {
    $ws.Range("H132").Value = 3465.25
    $ws.Range("I132").Value = 3044.6
    $ws.Range("K132").Value = 9133.799999999999
    $ws.Range("M132").Value = -6603.799999999999
}

$ws = $wb.Worksheets.Item("WVR")
# Row 18: Leve Item ID (G18) = 3543
if ($ws.Cells.Item(18, 7).Value2 -eq 3543) {
    $ws.Range("H18").Value = 0
    $ws.Range("J18").Value = 0
    $ws.Range("L18").Value = 0
    $ws.Range("N18").ClearContents()
}
# Row 20: Leve Item ID (G20) = 3023
if ($ws.Cells.Item(20, 7).Value2 -eq 3023) {
    $ws.Range("H20").Value = 11899.5
    $ws.Range("J20").Value = 11899.5
    $ws.Range("L20").Value = 11899.5
    $ws.Range("N20").Value = -12379.5
}
# Row 107: Leve Item ID (G107) = 27746
if ($ws.Cells.Item(107, 7).Value2 -eq 27746) {
    $ws.Range("H107").Value = 652.5278
    $ws.Range("I107").Value = 622.3570999999999
    $ws.Range("J107").Value = 758.125
    $ws.Range("K107").Value = 1867.0713
    $ws.Range("L107").Value = 2274.375
    $ws.Range("M107").Value = 52.92870000000016
    $ws.Range("N107").Value = -6114.375
}
# Row 122: Leve Item ID (G122) = 36208
if ($ws.Cells.Item(122, 7).Value2 -eq 36208) {
    $ws.Range("H122").Value = 209044.02
    $ws.Range("I122").Value = 282126.66
    $ws.Range("K122").Value = 846379.98
    $ws.Range("M122").Value = -843929.98
}
# Row 126: Leve Item ID (G126) = 36210
if ($ws.Cells.Item(126, 7).Value2 -eq 36210) {
    $ws.Range("H126").Value = 2817.8333
    $ws.Range("I126").Value = 2216.5
    $ws.Range("J126").Value = 3419.1667
    $ws.Range("K126").Value = 6649.5
    $ws.Range("L126").Value = 10257.5001
    $ws.Range("M126").Value = -4179.5
    $ws.Range("N126").Value = -15197.5001
}

